$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the current "last row" date format (YYYY-MM-DD) before we touch anything
$lastRowDateFormat = $ws.Range("A7").NumberFormat

# Row 7 (2021-11-16) is no longer the last row, so it should now use the
# regular date format used by the rest of the date column.
$ws.Range("A7").NumberFormat = $ws.Range("A6").NumberFormat

# Append the new data row for 2021-11-17
$ws.Range("A8").Value = 44517
$ws.Range("B8").Value = -1973.7

# The newly appended row becomes the new "last row" and gets the distinct format
$ws.Range("A8").NumberFormat = $lastRowDateFormat
